$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with revised figures
$ws.Range("B173").Value = 2547
$ws.Range("C173").Value = 10556
$ws.Range("B174").Value = 2100
$ws.Range("B175").Value = 2913

# Append new row 176 with new period (01-07-2021 stored as text, not a date)
$ws.Range("A176").NumberFormat = "@"
$ws.Range("A176").Value = "01-07-2021"
$ws.Range("A176").Style = $ws.Range("A175").Style
$ws.Range("B176").Value = 3275
$ws.Range("C176").Value = 7657
